$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values would otherwise be
# auto-coerced to numbers by Excel (losing the exact string formatting).
$textCells = @("D4","D5","D6","D10","D11","D12","D17","D19","D22","D27","D29","D30","D32","D36","D37","D40","D42","D43","D45","D46","D47","D48","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "61.739.98"
$ws.Range("E2").Value = "  +2.97%  "
$ws.Range("D3").Value = "3.402.67"
$ws.Range("E3").Value = "  +4.02%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "576.14"
$ws.Range("E5").Value = "  +3.74%  "
$ws.Range("D6").Value = "138.59"
$ws.Range("E6").Value = "  +9.16%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "3.402.20"
$ws.Range("E8").Value = "  +3.87%  "
$ws.Range("E9").Value = "  +1.57%  "
$ws.Range("D10").Value = "7.52"
$ws.Range("E10").Value = "  +3.56%  "
$ws.Range("D11").Value = "0.127"
$ws.Range("E11").Value = "  +9.19%  "
$ws.Range("D12").Value = "0.396"
$ws.Range("E12").Value = "  +7.54%  "
$ws.Range("D13").Value = "3.972.81"
$ws.Range("E13").Value = "  +3.57%  "
$ws.Range("E14").Value = "  +1.73%  "
$ws.Range("E15").Value = "  +8.35%  "
$ws.Range("D16").Value = "3.393.16"
$ws.Range("E16").Value = "  +3.43%  "
$ws.Range("D17").Value = "25.39"
$ws.Range("E17").Value = "  +5.65%  "
$ws.Range("D18").Value = "61.755.30"
$ws.Range("E18").Value = "  +2.62%  "
$ws.Range("D19").Value = "14.08"
$ws.Range("E19").Value = "  +6.74%  "
$ws.Range("E20").Value = "  +5.93%  "
$ws.Range("E21").Value = "  +4.60%  "
$ws.Range("D22").Value = "389.02"
$ws.Range("E22").Value = "  +11.39%  "
$ws.Range("E23").Value = "  +4.67%  "
$ws.Range("D24").Value = "3.531.69"
$ws.Range("E24").Value = "  +3.86%  "
$ws.Range("E25").Value = "  +18.61%  "
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("D27").Value = "71.17"
$ws.Range("E27").Value = "  +3.02%  "
$ws.Range("E28").Value = "  +14.68%  "
$ws.Range("D29").Value = "7.75"
$ws.Range("E29").Value = "  +7.75%  "
$ws.Range("D30").Value = "0.998"
$ws.Range("E30").Value = "  -0.26%  "
$ws.Range("E31").Value = "  +7.81%  "
$ws.Range("D32").Value = "0.159"
$ws.Range("E32").Value = "  +6.68%  "
$ws.Range("E33").Value = "  +4.04%  "
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("D35").Value = "3.426.69"
$ws.Range("E35").Value = "  +3.69%  "
$ws.Range("D36").Value = "23.55"
$ws.Range("E36").Value = "  +4.59%  "
$ws.Range("D37").Value = "5.52"
$ws.Range("E37").Value = "  +4.73%  "
$ws.Range("E38").Value = "  +4.18%  "
$ws.Range("E39").Value = "  +5.81%  "
$ws.Range("D40").Value = "162.18"
$ws.Range("E40").Value = "  +3.01%  "
$ws.Range("E41").Value = "  +7.17%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "1.75"
$ws.Range("E42").Value = "  +14.19%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  -0.29%  "
$ws.Range("E44").Value = "  +8.83%  "
$ws.Range("D45").Value = "0.774"
$ws.Range("E45").Value = "  +5.64%  "
$ws.Range("D46").Value = "4.45"
$ws.Range("E46").Value = "  +3.34%  "
$ws.Range("D47").Value = "41.39"
$ws.Range("E47").Value = "  +1.31%  "
$ws.Range("D48").Value = "24.71"
$ws.Range("E48").Value = "  +10.43%  "
$ws.Range("E49").Value = "  +5.82%  "
$ws.Range("D50").Value = "23.11"
$ws.Range("E50").Value = "  +9.18%  "
$ws.Range("D51").Value = "2.377.70"
$ws.Range("E51").Value = "  +10.61%  "
